# New weekly price report: insert a new record for "Locoto" (Vega Modelo de
# Temuco) right above the existing row 61, pushing the rest of the table
# down by one row (mirrors Excel's native "insert row" behaviour, which
# shifts all subsequent rows/styles down automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 61; everything from 61..112 moves to 62..113.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row with the latest observation.
$ws.Cells.Item(61, 1).Value  = 10
$ws.Cells.Item(61, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(61, 3).Value  = "La Araucanía"
$ws.Cells.Item(61, 4).Value  = 45240
$ws.Cells.Item(61, 5).Value  = 9
$ws.Cells.Item(61, 6).Value  = 100112042
$ws.Cells.Item(61, 7).Value  = "Locoto"
$ws.Cells.Item(61, 8).Value  = "Sin especificar"
$ws.Cells.Item(61, 9).Value  = "Primera"
$ws.Cells.Item(61, 10).Value = 100
$ws.Cells.Item(61, 11).Value = 3800
$ws.Cells.Item(61, 12).Value = 3800
$ws.Cells.Item(61, 13).Value = 3800
$ws.Cells.Item(61, 14).Value = "$/kilo"
$ws.Cells.Item(61, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(61, 16).Value = 3800
$ws.Cells.Item(61, 17).Value = 1
$ws.Cells.Item(61, 18).Value = "Hortaliza"
